$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 143
$ws.Cells.Item(143, 2).Value = 7532414
$ws.Cells.Item(143, 5).Value = 'Independiente Petrolero'
$ws.Cells.Item(143, 6).Value = 'Real Santa Cruz'
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 0
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 11).Value = 'H'
$ws.Cells.Item(143, 12).Value = 1.571
$ws.Cells.Item(143, 13).Value = 3.75
$ws.Cells.Item(143, 14).Value = 5
$ws.Cells.Item(143, 15).Value = 1.3
$ws.Cells.Item(143, 16).Value = 5
$ws.Cells.Item(143, 17).Value = 11
$ws.Cells.Item(143, 18).Value = -1.75
$ws.Cells.Item(143, 19).Value = 2
$ws.Cells.Item(143, 20).Value = 1.8
$ws.Cells.Item(143, 21).Value = 3
$ws.Cells.Item(143, 22).Value = 1.85
$ws.Cells.Item(143, 23).Value = 1.95
$ws.Cells.Item(143, 24).Value = 0.3
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 28).Value = 0.8
$ws.Cells.Item(143, 30).Value = 0.95

# Row 144
$ws.Cells.Item(144, 2).Value = 7532413
$ws.Cells.Item(144, 5).Value = 'Libertad Gran Mamore FC'
$ws.Cells.Item(144, 6).Value = 'Club Aurora'
$ws.Cells.Item(144, 8).Value = 1
$ws.Cells.Item(144, 10).Value = 0
$ws.Cells.Item(144, 12).Value = 2.25
$ws.Cells.Item(144, 13).Value = 3.3
$ws.Cells.Item(144, 14).Value = 2.8
$ws.Cells.Item(144, 15).Value = 2.375
$ws.Cells.Item(144, 16).Value = 3.4
$ws.Cells.Item(144, 19).Value = 2.025
$ws.Cells.Item(144, 20).Value = 1.775
$ws.Cells.Item(144, 21).Value = 2.5
$ws.Cells.Item(144, 22).Value = 1.9
$ws.Cells.Item(144, 23).Value = 1.9
$ws.Cells.Item(144, 28).Value = 0.7749999999999999
$ws.Cells.Item(144, 29).Value = -1
$ws.Cells.Item(144, 30).Value = 0.8999999999999999

# Row 145
$ws.Cells.Item(145, 2).Value = 7532412
$ws.Cells.Item(145, 5).Value = 'Vaca Diez'
$ws.Cells.Item(145, 6).Value = 'Blooming'
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 3
$ws.Cells.Item(145, 9).Value = 0
$ws.Cells.Item(145, 10).Value = 2
$ws.Cells.Item(145, 11).Value = 'A'
$ws.Cells.Item(145, 12).Value = 1.727
$ws.Cells.Item(145, 14).Value = 4
$ws.Cells.Item(145, 15).Value = 2.3
$ws.Cells.Item(145, 16).Value = 3.6
$ws.Cells.Item(145, 17).Value = 2.875
$ws.Cells.Item(145, 18).Value = -0.25
$ws.Cells.Item(145, 19).Value = 1.95
$ws.Cells.Item(145, 20).Value = 1.85
$ws.Cells.Item(145, 21).Value = 2.75
$ws.Cells.Item(145, 22).Value = 1.925
$ws.Cells.Item(145, 23).Value = 1.875
$ws.Cells.Item(145, 24).Value = -1
$ws.Cells.Item(145, 26).Value = 1.875
$ws.Cells.Item(145, 28).Value = 0.8500000000000001
$ws.Cells.Item(145, 29).Value = 0.4625
$ws.Cells.Item(145, 30).Value = -0.5

# Row 148
$ws.Cells.Item(148, 2).Value = 7532419
$ws.Cells.Item(148, 5).Value = 'Oriente Petrolero'
$ws.Cells.Item(148, 6).Value = 'Jorge Wilstermann'
$ws.Cells.Item(148, 12).Value = 2.2
$ws.Cells.Item(148, 13).Value = 2.5
$ws.Cells.Item(148, 14).Value = 4.5
$ws.Cells.Item(148, 15).Value = 2.375
$ws.Cells.Item(148, 16).Value = 2.45
$ws.Cells.Item(148, 17).Value = 4.5
$ws.Cells.Item(148, 18).Value = -0.25
$ws.Cells.Item(148, 19).Value = 1.9
$ws.Cells.Item(148, 20).Value = 1.9
$ws.Cells.Item(148, 21).Value = 2
$ws.Cells.Item(148, 24).Value = 1.375
$ws.Cells.Item(148, 27).Value = 0.8999999999999999
$ws.Cells.Item(148, 29).Value = 0.95
$ws.Cells.Item(148, 30).Value = -1

# Row 149
$ws.Cells.Item(149, 2).Value = 7532420
$ws.Cells.Item(149, 5).Value = 'Club Aurora'
$ws.Cells.Item(149, 6).Value = 'Vaca Diez'
$ws.Cells.Item(149, 7).Value = 3
$ws.Cells.Item(149, 9).Value = 2
$ws.Cells.Item(149, 12).Value = 1.333
$ws.Cells.Item(149, 13).Value = 5
$ws.Cells.Item(149, 14).Value = 8
$ws.Cells.Item(149, 15).Value = 1.3
$ws.Cells.Item(149, 16).Value = 6.5
$ws.Cells.Item(149, 17).Value = 7
$ws.Cells.Item(149, 19).Value = 1.8
$ws.Cells.Item(149, 20).Value = 2
$ws.Cells.Item(149, 21).Value = 3.25
$ws.Cells.Item(149, 22).Value = 1.95
$ws.Cells.Item(149, 23).Value = 1.85
$ws.Cells.Item(149, 24).Value = 0.3
$ws.Cells.Item(149, 27).Value = 0.8
$ws.Cells.Item(149, 29).Value = -0.5
$ws.Cells.Item(149, 30).Value = 0.425

# Row 150
$ws.Cells.Item(150, 2).Value = 7532421
$ws.Cells.Item(150, 5).Value = 'Guabira'
$ws.Cells.Item(150, 6).Value = 'Independiente Petrolero'
$ws.Cells.Item(150, 7).Value = 2
$ws.Cells.Item(150, 9).Value = 1
$ws.Cells.Item(150, 12).Value = 1.4
$ws.Cells.Item(150, 13).Value = 4.5
$ws.Cells.Item(150, 14).Value = 7.5
$ws.Cells.Item(150, 15).Value = 1.333
$ws.Cells.Item(150, 16).Value = 5.5
$ws.Cells.Item(150, 17).Value = 9.5
$ws.Cells.Item(150, 18).Value = -1.5
$ws.Cells.Item(150, 19).Value = 1.85
$ws.Cells.Item(150, 20).Value = 1.95
$ws.Cells.Item(150, 21).Value = 3
$ws.Cells.Item(150, 22).Value = 1.825
$ws.Cells.Item(150, 23).Value = 1.975
$ws.Cells.Item(150, 24).Value = 0.333
$ws.Cells.Item(150, 27).Value = 0.8500000000000001
$ws.Cells.Item(150, 29).Value = -1
$ws.Cells.Item(150, 30).Value = 0.9750000000000001

# Row 214
$ws.Cells.Item(214, 2).Value = 8038943
$ws.Cells.Item(214, 5).Value = 'San Jose de Oruro'
$ws.Cells.Item(214, 6).Value = 'Bolivar'
$ws.Cells.Item(214, 9).Value = 1
$ws.Cells.Item(214, 12).Value = 2.3
$ws.Cells.Item(214, 13).Value = 3.5
$ws.Cells.Item(214, 14).Value = 2.625
$ws.Cells.Item(214, 15).Value = 2.8
$ws.Cells.Item(214, 16).Value = 3.6
$ws.Cells.Item(214, 17).Value = 2.375
$ws.Cells.Item(214, 18).Value = 0.25
$ws.Cells.Item(214, 19).Value = 1.8
$ws.Cells.Item(214, 20).Value = 2
$ws.Cells.Item(214, 21).Value = 3.25
$ws.Cells.Item(214, 22).Value = 1.975
$ws.Cells.Item(214, 23).Value = 1.825
$ws.Cells.Item(214, 24).Value = 1.8
$ws.Cells.Item(214, 27).Value = 0.8
$ws.Cells.Item(214, 28).Value = -1
$ws.Cells.Item(214, 29).Value = -0.5
$ws.Cells.Item(214, 30).Value = 0.4125

# Row 215
$ws.Cells.Item(215, 2).Value = 8039392
$ws.Cells.Item(215, 5).Value = 'Oriente Petrolero'
$ws.Cells.Item(215, 6).Value = 'Jorge Wilstermann'
$ws.Cells.Item(215, 9).Value = 0
$ws.Cells.Item(215, 12).Value = 2
$ws.Cells.Item(215, 13).Value = 3.25
$ws.Cells.Item(215, 14).Value = 3.4
$ws.Cells.Item(215, 15).Value = 1.727
$ws.Cells.Item(215, 16).Value = 4
$ws.Cells.Item(215, 17).Value = 4.5
$ws.Cells.Item(215, 18).Value = -0.75
$ws.Cells.Item(215, 19).Value = 1.9
$ws.Cells.Item(215, 20).Value = 1.9
$ws.Cells.Item(215, 21).Value = 2.75
$ws.Cells.Item(215, 22).Value = 1.9
$ws.Cells.Item(215, 23).Value = 1.9
$ws.Cells.Item(215, 24).Value = 0.7270000000000001
$ws.Cells.Item(215, 27).Value = 0.45
$ws.Cells.Item(215, 28).Value = -0.5
$ws.Cells.Item(215, 29).Value = 0.45
$ws.Cells.Item(215, 30).Value = -0.5

# Row 262
$ws.Cells.Item(262, 15).Value = 2.35
$ws.Cells.Item(262, 17).Value = 2.9
$ws.Cells.Item(262, 19).Value = 1.7
$ws.Cells.Item(262, 20).Value = 2.1
